$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Answer")
$ws.Range("C1").Value = "Mark"
$ws.Range("C2").Select()
